$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 2 (first data row) - remaining rows shift up by one
$ws.Rows.Item(2).Delete()

# Append new gyroscope data rows at the bottom (new rows 21-31)
$newData = @(
    @(0.06544841687701197, -3.185547882748625, -1.701872078413778),
    @(-8.064484859250227, -7.12666956174002, 2.896797829067585),
    @(-3.302934056704871, -2.102355150832409, -0.8136555189938615),
    @(3.127995854800499, -0.3605925717305229, -0.3430995842844817),
    @(5.373623056510048, -6.720149974233094, -0.7005343388036227),
    @(3.014144754901412, -5.344652057923007, -1.741643652473547),
    @(-2.249343609072485, -2.945478901420671, 1.552763677134959),
    @(-3.602893884648981, -3.38362657900938, -0.5896314640635061),
    @(-4.196339511379769, -2.290577345287686, 1.337645951005574),
    @(0.4132739001328325, -1.036843425527086, 6.001093726797151),
    @(3.482481982290137, 3.476468096074395, 0.00721320909322376)
)

$startRow = 21
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
